$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Price" (column D) values, forcing text storage so values
# like "1.241" or "0.3833" are not reinterpreted as numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '23.531.63'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.655.57'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '302.35'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3833'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '51.15'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.241'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08187'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.39'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.486'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.486'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.656.46'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '97.42'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06970'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.826'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.67'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '23.545.43'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.500'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.991'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.18'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '152.67'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.238'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.72'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.841.05'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.190'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.252'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '12.14'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.054'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02798'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.122'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2497'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.08790'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.07016'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.19'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7013'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.336'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.09'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6542'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07917'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '127.91'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.192'
$ws.Range('D51').Style = 'Normal'

# Update "Volume(1h)" (column E) percentage-change values.
$ws.Range('E2').Value = '  +1.50%  '
$ws.Range('E3').Value = '  +2.70%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('E7').Value = '  +1.35%  '
$ws.Range('E8').Value = '  -1.67%  '
$ws.Range('E9').Value = '  +1.70%  '
$ws.Range('E10').Value = '  +3.09%  '
$ws.Range('E11').Value = '  +1.10%  '
$ws.Range('E12').Value = '  -0.24%  '
$ws.Range('E13').Value = '  +1.25%  '
$ws.Range('E14').Value = '  +2.11%  '
$ws.Range('E15').Value = '  +3.19%  '
$ws.Range('E16').Value = '  +0.99%  '
$ws.Range('E17').Value = '  +3.60%  '
$ws.Range('E18').Value = '  +3.32%  '
$ws.Range('E19').Value = '  +0.90%  '
$ws.Range('E20').Value = '  +5.18%  '
$ws.Range('E21').Value = '  +2.84%  '
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('E23').Value = '  +3.10%  '
$ws.Range('E24').Value = '  +1.60%  '
$ws.Range('E25').Value = '  -0.23%  '
$ws.Range('E26').Value = '  -0.54%  '
$ws.Range('E27').Value = '  +1.51%  '
$ws.Range('E28').Value = '  +1.20%  '
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('E30').Value = '  +1.20%  '
$ws.Range('E31').Value = '  +3.28%  '
$ws.Range('E32').Value = '  +11.23%  '
$ws.Range('E33').Value = '  +7.27%  '
$ws.Range('E34').Value = '  +7.49%  '
$ws.Range('E35').Value = '  -0.99%  '
$ws.Range('E36').Value = '  +3.60%  '
$ws.Range('E37').Value = '  +4.84%  '
$ws.Range('E38').Value = '  +1.90%  '
$ws.Range('E39').Value = '  +0.61%  '
$ws.Range('E40').Value = '  +1.34%  '
$ws.Range('E41').Value = '  +10.64%  '
$ws.Range('E42').Value = '  +2.01%  '
$ws.Range('E43').Value = '  +1.07%  '
$ws.Range('E44').Value = '  +5.37%  '
$ws.Range('E45').Value = '  +3.71%  '
$ws.Range('E46').Value = '  -0.10%  '
$ws.Range('E47').Value = '  +2.76%  '
$ws.Range('E48').Value = '  +0.26%  '
$ws.Range('E49').Value = '  +0.72%  '
$ws.Range('E50').Value = '  +1.25%  '
$ws.Range('E51').Value = '  +2.02%  '
